$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 69

# Column A holds a date-like string that must stay plain text (matching the
# existing rows), not get auto-converted into a date serial number. Prefixing
# with an apostrophe forces text entry (same as typing it in the Excel UI);
# resetting the Style afterwards drops the "quote prefix" cell format that
# Excel tacks on, so the cell ends up with no explicit style - just like the
# other data rows.
$ws.Cells.Item($row, 1).Value = "'2025-10-23"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 54.52999877929688
$ws.Cells.Item($row, 3).Value = 405.8500061035156
$ws.Cells.Item($row, 4).Value = 328.3500061035156
